# Insert a new "ePhone" column between "Tel" (N) and "Email" (O) on the
# company sheet, pushing Email / Industry / Tax Duration one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing "Email" column (and everything after it) one column to the
# right by inserting a blank column at O, then fill in the new header.
$ws.Range("O1").EntireColumn.Insert()
$ws.Range("O1").Value = "ePhone"
